# Bekir Tanriover_2023.xlsx data correction
# - cited_by_count (M) bumped for rows 2 and 3 (25->26, 5->6)
# - new works added: dept1703 (Marek Rychlik / W4386348049) and
#   dept1704 (Marek Rychlik / W4386396309); existing rows shift down
#   so row 8 (Azhar / W4386574113) becomes row 10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# cited_by_count updates
Set-TextValue $ws.Range("M2") "26"
Set-TextValue $ws.Range("M3") "6"

# Row 8 <- dept1703: Marek Rychlik et al., W4386348049 (arXiv preprint)
$ws.Range("A8").Value = "Marek Rychlik, Bekir Tanrıöver, Yuxing Han"
$ws.Range("B8").Value = "; ; "
$ws.Range("C8").Value = "https://openalex.org/W4386348049"
$ws.Range("D8").Value = "Large-scale data extraction from the UNOS organ donor documents"
Set-TextValue $ws.Range("E8") "2023-08-30"
$ws.Range("F8").Value = "arXiv (Cornell University)"
$ws.Range("G8").Value = "Cornell University"
$ws.Range("H8").Value = "https://arxiv.org/abs/2308.15752"
$ws.Range("I8").Value = "N/A"
$ws.Range("J8").Value = "submittedVersion"
$ws.Range("K8").Value = "green"
$ws.Range("L8").Value = "en"
Set-TextValue $ws.Range("M8") "0"
Set-TextValue $ws.Range("N8") "2023"
$ws.Range("O8").Value = "NA"
$ws.Range("P8").Value = "https://doi.org/10.48550/arxiv.2308.15752"
$ws.Range("Q8").Value = "article"

# Row 9 <- dept1704: Marek Rychlik et al., W4386396309 (Preprints.org)
$ws.Range("A9").Value = "Marek Rychlik, Bekir Tanrıöver, Yuxing Han"
$ws.Range("B9").Value = "; ; "
$ws.Range("C9").Value = "https://openalex.org/W4386396309"
$ws.Range("D9").Value = "Large-Scale Data Extraction from the UNOS Organ Donor Documents"
Set-TextValue $ws.Range("E9") "2023-08-31"
$ws.Range("F9").Value = "N/A"
$ws.Range("G9").Value = "N/A"
$ws.Range("H9").Value = "https://doi.org/10.20944/preprints202308.2121.v1"
$ws.Range("I9").Value = "N/A"
$ws.Range("J9").Value = "submittedVersion"
$ws.Range("K9").Value = "bronze"
$ws.Range("L9").Value = "en"
Set-TextValue $ws.Range("M9") "0"
Set-TextValue $ws.Range("N9") "2023"
$ws.Range("O9").Value = "NA"
$ws.Range("P9").Value = "https://doi.org/10.20944/preprints202308.2121.v1"
$ws.Range("Q9").Value = "article"

# Row 10 <- Azhar et al., W4386574113 (shifted down from old row 8)
$ws.Range("A10").Value = "Ambreen Azhar, Bekir Tanrıöver, Ahmet B. Gungor, Miklós Molnár, Gaurav Gupta"
$ws.Range("B10").Value = "Department of Internal Medicine, Division of Nephrology, Virginia Commonwealth University, 1101 East Marshall Street, PO Box 980160, Richmond, VA, 23298, USA; Division of Nephrology, College of Medicine, University of Arizona, Tucson, AZ, USA; Division of Nephrology, College of Medicine, University of Arizona, Tucson, AZ, USA; Department of Internal Medicine, Division of Nephrology & Hypertension, University of Utah Spencer Fox Eccles School of Medicine, Salt Lake City, UT, USA; Department of Internal Medicine, Division of Nephrology, Virginia Commonwealth University, 1101 East Marshall Street, PO Box 980160, Richmond, VA, 23298, USA"
$ws.Range("C10").Value = "https://openalex.org/W4386574113"
$ws.Range("D10").Value = "Virologic Studies in COVID-Positive Donors"
Set-TextValue $ws.Range("E10") "2023-09-09"
$ws.Range("F10").Value = "Current Transplantation Reports"
$ws.Range("G10").Value = "Springer Science+Business Media"
$ws.Range("H10").Value = "https://doi.org/10.1007/s40472-023-00411-7"
$ws.Range("I10").Value = "N/A"
$ws.Range("J10").Value = "N/A"
$ws.Range("K10").Value = "closed"
$ws.Range("L10").Value = "en"
Set-TextValue $ws.Range("M10") "0"
Set-TextValue $ws.Range("N10") "2023"
$ws.Range("O10").Value = "NA"
$ws.Range("P10").Value = "https://doi.org/10.1007/s40472-023-00411-7"
$ws.Range("Q10").Value = "article"

